$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Insert two new rows before row 12, pushing existing rows 12-18 down to 14-20.
$ws.Rows.Item(12).Resize(2).Insert()

# Fill in the two newly inserted rows with the new datatype fields.
$ws.Range("D12").Value = "Date"
$ws.Range("E12").Value = "foo"

$ws.Range("D13").Value = "Double"
$ws.Range("E13").Value = "bar"

# Update the SmartRules signature (now at row 16 after the insert) to mention the new field.
$ws.Range("D16").Value = "SmartRules MyDatatype myRules( String myCode, MyDatatype myObj, Date foo)"

# Update the worksheet dimension/selection to match the new layout.
$ws.Range("E14").Select()
